{"js": "// Replace the worksheet date and every two-digit-by-two-digit multiplication\n// answer cell with the new values from the next day's generated sheet.\n// Each entry is [oldText, newText]; all old values are unique substrings in\n// the document, so a direct search+replace per pair is unambiguous.\nconst replacements = [\n  [\"2025-07-31 Thursday\", \"2025-08-01 Friday\"],\n  [\"58\u00d770=4060\", \"29\u00d779=2291\"],\n  [\"82\u00d761=5002\", \"21\u00d730=630\"],\n  [\"30\u00d777=2310\", \"74\u00d771=5254\"],\n  [\"93\u00d769=6417\", \"80\u00d773=5840\"],\n  [\"87\u00d739=3393\", \"11\u00d722=242\"],\n  [\"38\u00d788=3344\", \"59\u00d718=1062\"],\n  [\"50\u00d798=4900\", \"45\u00d784=3780\"],\n  [\"67\u00d738=2546\", \"27\u00d725=675\"],\n  [\"44\u00d742=1848\", \"13\u00d796=1248\"],\n  [\"37\u00d725=925\", \"44\u00d717=748\"],\n  [\"27\u00d797=2619\", \"38\u00d776=2888\"],\n  [\"55\u00d758=3190\", \"43\u00d791=3913\"],\n  [\"38\u00d726=988\", \"21\u00d714=294\"],\n  [\"76\u00d737=2812\", \"54\u00d797=5238\"],\n  [\"20\u00d727=540\", \"73\u00d775=5475\"],\n  [\"88\u00d738=3344\", \"59\u00d787=5133\"],\n  [\"18\u00d759=1062\", \"20\u00d767=1340\"],\n  [\"34\u00d795=3230\", \"61\u00d767=4087\"],\n  [\"43\u00d757=2451\", \"37\u00d722=814\"],\n  [\"82\u00d747=3854\", \"26\u00d736=936\"],\n  [\"76\u00d785=6460\", \"39\u00d789=3471\"],\n  [\"71\u00d772=5112\", \"16\u00d787=1392\"],\n  [\"45\u00d798=4410\", \"27\u00d790=2430\"],\n  [\"99\u00d786=8514\", \"59\u00d728=1652\"],\n  [\"76\u00d794=7144\", \"84\u00d770=5880\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const rng of results.items) {\n    rng.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the worksheet date and every two-digit-by-two-digit multiplication\n# answer cell with the new values from the next day's generated sheet.\n# Each pair is (oldText, newText); all old values are unique in the document,\n# so a whole-document Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-07-31 Thursday\", \"2025-08-01 Friday\")\n    ,@(\"58\u00d770=4060\", \"29\u00d779=2291\")\n    ,@(\"82\u00d761=5002\", \"21\u00d730=630\")\n    ,@(\"30\u00d777=2310\", \"74\u00d771=5254\")\n    ,@(\"93\u00d769=6417\", \"80\u00d773=5840\")\n    ,@(\"87\u00d739=3393\", \"11\u00d722=242\")\n    ,@(\"38\u00d788=3344\", \"59\u00d718=1062\")\n    ,@(\"50\u00d798=4900\", \"45\u00d784=3780\")\n    ,@(\"67\u00d738=2546\", \"27\u00d725=675\")\n    ,@(\"44\u00d742=1848\", \"13\u00d796=1248\")\n    ,@(\"37\u00d725=925\", \"44\u00d717=748\")\n    ,@(\"27\u00d797=2619\", \"38\u00d776=2888\")\n    ,@(\"55\u00d758=3190\", \"43\u00d791=3913\")\n    ,@(\"38\u00d726=988\", \"21\u00d714=294\")\n    ,@(\"76\u00d737=2812\", \"54\u00d797=5238\")\n    ,@(\"20\u00d727=540\", \"73\u00d775=5475\")\n    ,@(\"88\u00d738=3344\", \"59\u00d787=5133\")\n    ,@(\"18\u00d759=1062\", \"20\u00d767=1340\")\n    ,@(\"34\u00d795=3230\", \"61\u00d767=4087\")\n    ,@(\"43\u00d757=2451\", \"37\u00d722=814\")\n    ,@(\"82\u00d747=3854\", \"26\u00d736=936\")\n    ,@(\"76\u00d785=6460\", \"39\u00d789=3471\")\n    ,@(\"71\u00d772=5112\", \"16\u00d787=1392\")\n    ,@(\"45\u00d798=4410\", \"27\u00d790=2430\")\n    ,@(\"99\u00d786=8514\", \"59\u00d728=1652\")\n    ,@(\"76\u00d794=7144\", \"84\u00d770=5880\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
